$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (pushes existing rows 15-25 down to 16-26)
$ws.Rows("15:15").Insert()

# Populate the new row 15 with the "Very Minor Works" line item
$ws.Range("A15").Value = "Very Minor Works"
$ws.Range("B15").Value = "p/year(s)"
$ws.Range("C15").Value = 1450

# Match number formatting of the row above before writing numeric values
$ws.Range("D15:H15").NumberFormat = $ws.Range("D14").NumberFormat
$ws.Range("D15").Value = 0
$ws.Range("E15").Formula = "=C15*D15"
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Formula = "=E15-F15-G15"
$ws.Range("M15").Value = 1

# Widen column A (bestfit-style width for the new, longer label)
$ws.Columns("A:A").ColumnWidth = 24.83

# Data validation on B14 now also covers the new B15 cell
$ws.Range("B14:B15").Select()

# Final selection used by the author when saving the workbook
$ws.Range("M15").Select()
